# fix broken mapping_orig.txt links
# Change the R column ("mapping-orig.txt" link) formulas on the
# dataset_metadata sheet from the old GitHub-pages relative link
# ("./datasets/...") to a relative link that resolves against
# github.com/.../blob/master/... ("../../blob/master/datasets/...").
# Also update the view state: dataset_metadata becomes the active /
# selected sheet (instead of "tasks"), with its selection moved to R17.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("dataset_metadata")

for ($r = 2; $r -le 16; $r++) {
    $formula = '=CONCATENATE("../../blob/master/datasets/", B' + $r + ', "/mapping-orig.txt")'
    $wsData.Range("R$r").Formula = $formula
}

# Make dataset_metadata the active sheet (this also clears tabSelected
# on whichever sheet was previously active, e.g. "tasks") and move the
# selection to R17, matching the new view state.
$wsData.Activate() | Out-Null
$wsData.Range("R17").Select() | Out-Null
